$p = $ppt.ActivePresentation

# --- Slide master: update footer / slide-number / date placeholders to show ---
# --- the generic template tokens instead of blank text / a stale literal "16" ---
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame -eq -1 -and $sh.PlaceholderFormat.Type -ne 0) {
        switch ($sh.PlaceholderFormat.Type) {
            15 { $sh.TextFrame.TextRange.Text = "<footer>" }      # ppPlaceholderFooter
            13 { $sh.TextFrame.TextRange.Text = "<number>" }      # ppPlaceholderSlideNumber
            16 { $sh.TextFrame.TextRange.Text = "<date/time>" }   # ppPlaceholderDate
        }
    }
}

# --- Slide 1: nudge the small Fibonacci logo picture next to the footer text ---
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Type -eq 13) {
        $sh.Left = 234.0
        $sh.Top = 373.2944881889764
        $sh.Width = 22.70552
        $sh.Height = 22.70552
    }
}
